$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "slide\ slick-slide\ slick-active"
$ws.Range("B3").Value = "slide\ slick-slide\ slick-active"
$ws.Range("C3").Value = "slide\ slick-slide\ slick-active"
$ws.Range("D3").Value = 'slick-list\ draggable"] [class="slide\ slick-slide'
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("I3").Value = " "
$ws.Range("J3").Value = " "
$ws.Range("K3:N3").NumberFormat = "@"
$ws.Range("K3").Value = "2"
$ws.Range("L3").Value = "4"
$ws.Range("M3").Value = "3"
$ws.Range("N3").Value = "4"
$ws.Range("O3").Value = " "
$ws.Range("P3").Value = "more\ slick-active"
$ws.Range("Q3").Value = " "
